$wb = $excel.ActiveWorkbook

# --- Sheet "Schedule": update computed Cost ($) and Unit Cost ($/ML) columns ---
$ws1 = $wb.Worksheets.Item("Schedule")
$ws1.Cells.Item(2, 5).Value2 = 335.0231820000001
$ws1.Cells.Item(2, 6).Value2 = 7.385872619047622
$ws1.Cells.Item(3, 5).Value2 = 409.49181975
$ws1.Cells.Item(3, 6).Value2 = 27.08279231150794

# --- Sheet "Detailed": refresh DateTime/Price/Pump_Status values (run 157) ---
$ws2 = $wb.Worksheets.Item("Detailed")
$ws2.Cells.Item(2, 1).Value2 = 46044.02083333334
$ws2.Cells.Item(2, 2).Value2 = 75.30216
$ws2.Cells.Item(3, 1).Value2 = 46044.04166666666
$ws2.Cells.Item(3, 2).Value2 = 73.2
$ws2.Cells.Item(4, 1).Value2 = 46044.0625
$ws2.Cells.Item(4, 2).Value2 = 77.94
$ws2.Cells.Item(5, 1).Value2 = 46044.08333333334
$ws2.Cells.Item(5, 2).Value2 = 73.2
$ws2.Cells.Item(6, 1).Value2 = 46044.10416666666
$ws2.Cells.Item(6, 2).Value2 = 73.2
$ws2.Cells.Item(7, 1).Value2 = 46044.125
$ws2.Cells.Item(7, 2).Value2 = 66.53829
$ws2.Cells.Item(8, 1).Value2 = 46044.14583333334
$ws2.Cells.Item(8, 2).Value2 = 73.2
$ws2.Cells.Item(9, 1).Value2 = 46044.16666666666
$ws2.Cells.Item(9, 2).Value2 = 73.2
$ws2.Cells.Item(9, 5).Value = "ON"
$ws2.Cells.Item(10, 1).Value2 = 46044.1875
$ws2.Cells.Item(10, 2).Value2 = 77.44047999999999
$ws2.Cells.Item(11, 1).Value2 = 46044.20833333334
$ws2.Cells.Item(11, 2).Value2 = 67.61373
$ws2.Cells.Item(12, 1).Value2 = 46044.22916666666
$ws2.Cells.Item(12, 2).Value2 = 78
$ws2.Cells.Item(13, 1).Value2 = 46044.25
$ws2.Cells.Item(13, 2).Value2 = 73.2
$ws2.Cells.Item(14, 1).Value2 = 46044.27083333334
$ws2.Cells.Item(14, 2).Value2 = 57.06
$ws2.Cells.Item(15, 1).Value2 = 46044.29166666666
$ws2.Cells.Item(15, 2).Value2 = 25.85393
$ws2.Cells.Item(16, 1).Value2 = 46044.3125
$ws2.Cells.Item(16, 2).Value2 = 0
$ws2.Cells.Item(17, 1).Value2 = 46044.33333333334
$ws2.Cells.Item(17, 2).Value2 = -5.74313
$ws2.Cells.Item(18, 1).Value2 = 46044.35416666666
$ws2.Cells.Item(18, 2).Value2 = -6.36915
$ws2.Cells.Item(19, 1).Value2 = 46044.375
$ws2.Cells.Item(19, 2).Value2 = -7.24179
$ws2.Cells.Item(20, 1).Value2 = 46044.39583333334
$ws2.Cells.Item(20, 2).Value2 = -7.7397
$ws2.Cells.Item(21, 1).Value2 = 46044.41666666666
$ws2.Cells.Item(21, 2).Value2 = -7.61678
$ws2.Cells.Item(22, 1).Value2 = 46044.4375
$ws2.Cells.Item(22, 2).Value2 = -8.976139999999999
$ws2.Cells.Item(23, 1).Value2 = 46044.45833333334
$ws2.Cells.Item(23, 2).Value2 = -8.919919999999999
$ws2.Cells.Item(24, 1).Value2 = 46044.47916666666
$ws2.Cells.Item(24, 2).Value2 = -5.58973
$ws2.Cells.Item(25, 1).Value2 = 46044.5
$ws2.Cells.Item(25, 2).Value2 = -8.02056
$ws2.Cells.Item(26, 1).Value2 = 46044.52083333334
$ws2.Cells.Item(26, 2).Value2 = -8.056279999999999
$ws2.Cells.Item(27, 1).Value2 = 46044.54166666666
$ws2.Cells.Item(27, 2).Value2 = -7.81277
$ws2.Cells.Item(28, 1).Value2 = 46044.5625
$ws2.Cells.Item(28, 2).Value2 = -6.01072
$ws2.Cells.Item(29, 1).Value2 = 46044.58333333334
$ws2.Cells.Item(29, 2).Value2 = -5.50985
$ws2.Cells.Item(30, 1).Value2 = 46044.60416666666
$ws2.Cells.Item(30, 2).Value2 = -3.6481
$ws2.Cells.Item(31, 1).Value2 = 46044.625
$ws2.Cells.Item(31, 2).Value2 = 0.51
$ws2.Cells.Item(32, 1).Value2 = 46044.64583333334
$ws2.Cells.Item(32, 2).Value2 = -12.01
$ws2.Cells.Item(33, 1).Value2 = 46044.66666666666
$ws2.Cells.Item(33, 2).Value2 = -10
$ws2.Cells.Item(33, 5).Value = "OFF"
$ws2.Cells.Item(34, 1).Value2 = 46044.6875
$ws2.Cells.Item(34, 2).Value2 = -12.01
$ws2.Cells.Item(35, 1).Value2 = 46044.70833333334
$ws2.Cells.Item(35, 2).Value2 = -11.01
$ws2.Cells.Item(36, 1).Value2 = 46044.72916666666
$ws2.Cells.Item(36, 2).Value2 = -8.107519999999999
$ws2.Cells.Item(37, 1).Value2 = 46044.75
$ws2.Cells.Item(37, 2).Value2 = -0.48574
$ws2.Cells.Item(38, 1).Value2 = 46044.77083333334
$ws2.Cells.Item(38, 2).Value2 = 3.05998
$ws2.Cells.Item(39, 1).Value2 = 46044.79166666666
$ws2.Cells.Item(39, 2).Value2 = 18.9562
$ws2.Cells.Item(40, 1).Value2 = 46044.8125
$ws2.Cells.Item(40, 2).Value2 = 55.33036
$ws2.Cells.Item(41, 1).Value2 = 46044.83333333334
$ws2.Cells.Item(41, 2).Value2 = 48.53134
$ws2.Cells.Item(41, 5).Value = "ON"
$ws2.Cells.Item(42, 1).Value2 = 46044.85416666666
$ws2.Cells.Item(42, 2).Value2 = 53.88034
$ws2.Cells.Item(43, 1).Value2 = 46044.875
$ws2.Cells.Item(43, 2).Value2 = 57.04367
$ws2.Cells.Item(44, 1).Value2 = 46044.89583333334
$ws2.Cells.Item(44, 2).Value2 = 57.03541
$ws2.Cells.Item(45, 1).Value2 = 46044.91666666666
$ws2.Cells.Item(45, 2).Value2 = 32.93261
$ws2.Cells.Item(46, 1).Value2 = 46044.9375
$ws2.Cells.Item(46, 2).Value2 = 56.98
$ws2.Cells.Item(47, 1).Value2 = 46044.95833333334
$ws2.Cells.Item(47, 2).Value2 = 56.98
$ws2.Cells.Item(48, 1).Value2 = 46044.97916666666
$ws2.Cells.Item(48, 2).Value2 = 56.60824

# Drop the now-obsolete last forecast row (the series shifted up by one period)
$ws2.Rows.Item(49).Delete()
